# Insert a new "setup" worksheet between "component_name" (sheet 4) and
# "heats" (currently sheet 5), populate it, select cell F8 on it, and make
# it the active sheet/tab (matches the workbook becoming activeTab=4 and
# the new sheet picking up tabSelected="1" while component_name loses it).

$wb = $excel.ActiveWorkbook

$componentNameSheet = $wb.Worksheets.Item("component_name")

$setupSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $componentNameSheet)
$setupSheet.Name = "setup"

# Write cells in this order so new shared-string entries are interned as
# "Calorimeter", "Initial volume", "ampoule" (matching the target file).
$setupSheet.Range("A1").Value = "Calorimeter"
$setupSheet.Range("A2").Value = "Initial volume"
$setupSheet.Range("B1").Value = "ampoule"
$setupSheet.Range("B2").Value = 1

$setupSheet.Range("F8").Select()
